$d = $word.ActiveDocument
$wNs = "xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'"

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2) | Out-Null
}

# Rewrite a whole "What we (don't) like" bullet paragraph in place, via its
# exact OOXML, so the leading empty <w:r/> artifact run that precedes the
# text run survives (a plain Find/Replace on these un-styled runs collapses
# that empty run away).
function Replace-Bullet($oldText, $newText) {
    foreach ($p in $d.Paragraphs) {
        if ($p.Range.Text -eq ($oldText + "`r")) {
            $xml = "<w:p $wNs><w:pPr><w:pStyle w:val='ListBullet'/><w:spacing w:line='240' w:lineRule='auto'/><w:ind w:left='720'/></w:pPr><w:r/><w:r><w:t>$newText</w:t></w:r></w:p>"
            $p.Range.InsertXML($xml) | Out-Null
            return $true
        }
    }
    return $false
}

# 1. Title: shorten / rewrite the H1 heading text.
Replace-Text "Play Fish Eye for Free: Review of Pragmatic Play's Unique Slot" "Play Fish Eye - Free Egypt-Based Slot Game"

# 2. Remove the whole "Meta description" paragraph (bold label + description text).
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "Meta description*") {
        $p.Range.Delete()
        break
    }
}

# 3 & 4. "What we like" bullets: reword first bullet, and move the old first
# bullet's wording down into the second bullet slot.
Replace-Bullet "Impressive graphics and detailed symbols" "Innovative take on the Egyptian theme"
Replace-Bullet "Innovative Egyptian-themed game set at the bottom of the Nile" "Impressive graphics and detailed symbols"

# 5. "What we don't like" bullets: reword both.
Replace-Bullet "Less frequent wins due to high volatility" "Less frequent wins on average"
Replace-Bullet "Not suitable for players who prefer low-volatility games" "Limited variety of special symbols"

# 6. Insert a new bold paragraph ("Play Fish Eye - Free Egypt-Based Slot Game")
# right after the last "What we don't like" bullet, before the italic
# image-prompt paragraph. Splitting the range at the end of the bullet's
# text (just before its paragraph mark) makes InsertXML add a genuine
# sibling paragraph instead of clobbering a neighbour.
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "Limited variety of special symbols*") {
        $pos = $p.Range.End - 1
        $r = $d.Range($pos, $pos)
        $r.InsertXML("<w:p $wNs><w:r/><w:r><w:rPr><w:b/></w:rPr><w:t>Play Fish Eye - Free Egypt-Based Slot Game</w:t></w:r></w:p>") | Out-Null
        break
    }
}

# 7. Rewrite the closing italic image-brief paragraph.
Replace-Text "Create a feature image that perfectly captures the excitement and adventure of Fish Eye by showcasing a happy Maya warrior wearing glasses. The image should feature the Maya warrior cheering as he dives into the depths of the Nile, surrounded by symbols from the game. The cartoon style should be colorful and engaging, drawing in players with its dynamic and fun design. The image could show the warrior holding a golden fish symbol, emphasizing the importance of these special symbols in the game. Overall, the feature image should convey the sense of exploration and discovery that players will experience when they dive into the world of Fish Eye." "Read our review of Fish Eye, a unique Egypt-based slot game to play for free."
